$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the entire contents of row 2 and row 3 (all columns A:AY).
$firstCol = 1
$lastCol = 51   # column AY

# Columns Y (25) and AA (27) hold plain "YYYY-MM-DD" text in this sheet;
# writing such a string back via Value2 would otherwise be auto-parsed into
# a date serial number. Mark just those cells as text first so the swapped
# values round-trip as plain strings, exactly like the source file.
$dateCols = @(25, 27)

for ($col = $firstCol; $col -le $lastCol; $col++) {
    $cell2 = $ws.Cells.Item(2, $col)
    $cell3 = $ws.Cells.Item(3, $col)

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    # Skip columns that are already blank on both rows -- writing an empty
    # string back would needlessly strip the (harmless) placeholder cell
    # that the original file stored there.
    if (($null -eq $v2) -and ($null -eq $v3)) {
        continue
    }

    if ($dateCols -contains $col) {
        $cell2.NumberFormat = "@"
        $cell3.NumberFormat = "@"
    }

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2

    if ($dateCols -contains $col) {
        $cell2.Style = "Normal"
        $cell3.Style = "Normal"
    }
}
